$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 114.65714
$ws.Range("I33").Value = 104.09677
$ws.Range("K33").Value = 104.09677
$ws.Range("M33").Value = 124.90323
$ws.Range("H55").Value = 170.66667
$ws.Range("I55").Value = 150
$ws.Range("K55").Value = 150
$ws.Range("M55").Value = 64
$ws.Range("H129").Value = 1001.3871
$ws.Range("I129").Value = 640.1429000000001
$ws.Range("J129").Value = 1106.75
$ws.Range("K129").Value = 1920.4287
$ws.Range("L129").Value = 3320.25
$ws.Range("M129").Value = 3079.5713
$ws.Range("N129").Value = -13320.25
$ws.Range("H135").Value = 659.1395
$ws.Range("I135").Value = 610.2895
$ws.Range("J135").Value = 1030.4
$ws.Range("K135").Value = 5492.6055
$ws.Range("L135").Value = 9273.6
$ws.Range("M135").Value = -2957.6055
$ws.Range("N135").Value = -14343.6
$ws.Range("H137").Value = 2027.091
$ws.Range("I137").Value = 1424.25
$ws.Range("K137").Value = 4272.75
$ws.Range("M137").Value = -1722.75
$ws.Range("H138").Value = 2503.012
$ws.Range("I138").Value = 1436.3469
$ws.Range("J138").Value = 4040.2646
$ws.Range("K138").Value = 4309.0407
$ws.Range("L138").Value = 12120.7938
$ws.Range("M138").Value = 830.9593000000004
$ws.Range("N138").Value = -22400.7938

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 41666.668
$ws.Range("J15").Value = 37500
$ws.Range("L15").Value = 37500
$ws.Range("N15").Value = -38200
$ws.Range("H32").Value = 10945.6875
$ws.Range("I32").Value = 11115.982
$ws.Range("K32").Value = 11115.982
$ws.Range("M32").Value = -10828.982
$ws.Range("H61").Value = 2712.25
$ws.Range("I61").Value = 1437.72
$ws.Range("K61").Value = 1437.72
$ws.Range("M61").Value = -1225.72
$ws.Range("H74").Value = 914.04443
$ws.Range("I74").Value = 1032.6774
$ws.Range("J74").Value = 651.3570999999999
$ws.Range("K74").Value = 1032.6774
$ws.Range("L74").Value = 651.3570999999999
$ws.Range("M74").Value = -158.6774
$ws.Range("N74").Value = -2399.3571
$ws.Range("H77").Value = 914.04443
$ws.Range("I77").Value = 1032.6774
$ws.Range("J77").Value = 651.3570999999999
$ws.Range("K77").Value = 5163.387000000001
$ws.Range("L77").Value = 3256.7855
$ws.Range("M77").Value = -795.3870000000006
$ws.Range("N77").Value = -11992.7855
$ws.Range("H97").Value = 1059
$ws.Range("I97").Value = 842.3
$ws.Range("J97").Value = 1781.3334
$ws.Range("K97").Value = 842.3
$ws.Range("L97").Value = 1781.3334
$ws.Range("M97").Value = -346.3
$ws.Range("N97").Value = -2773.3334
$ws.Range("H113").Value = 34199
$ws.Range("J113").Value = 34199
$ws.Range("L113").Value = 34199
$ws.Range("N113").Value = -42877
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H132").Value = 3242.5334
$ws.Range("I132").Value = 2892.375
$ws.Range("J132").Value = 3642.7144
$ws.Range("K132").Value = 8677.125
$ws.Range("L132").Value = 10928.1432
$ws.Range("M132").Value = -6147.125
$ws.Range("N132").Value = -15988.1432
$ws.Range("H136").Value = 2712.25
$ws.Range("I136").Value = 1437.72
$ws.Range("K136").Value = 4313.16
$ws.Range("M136").Value = -1763.16

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 85889.5
$ws.Range("J74").Value = 85889.5
$ws.Range("L74").Value = 85889.5
$ws.Range("N74").Value = -87761.5
$ws.Range("H77").Value = 85889.5
$ws.Range("J77").Value = 85889.5
$ws.Range("L77").Value = 257668.5
$ws.Range("N77").Value = -267028.5
$ws.Range("H94").Value = 41261.28
$ws.Range("I94").Value = 804.0769
$ws.Range("J94").Value = 85089.914
$ws.Range("K94").Value = 804.0769
$ws.Range("L94").Value = 85089.914
$ws.Range("M94").Value = -353.0769
$ws.Range("N94").Value = -85991.914
$ws.Range("H134").Value = 2911
$ws.Range("I134").Value = 2607.2727
$ws.Range("K134").Value = 7821.8181
$ws.Range("M134").Value = -5286.8181

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1717.4423
$ws.Range("I31").Value = 1338.9767
$ws.Range("K31").Value = 1338.9767
$ws.Range("M31").Value = -1043.9767
$ws.Range("H34").Value = 1717.4423
$ws.Range("I34").Value = 1338.9767
$ws.Range("K34").Value = 1338.9767
$ws.Range("M34").Value = -1136.9767
$ws.Range("H58").Value = 757208.0600000001
$ws.Range("I58").Value = 950749
$ws.Range("J58").Value = 2398.3
$ws.Range("K58").Value = 950749
$ws.Range("L58").Value = 2398.3
$ws.Range("M58").Value = -950546
$ws.Range("N58").Value = -2804.3
$ws.Range("H88").Value = 34417.5
$ws.Range("J88").Value = 34417.5
$ws.Range("L88").Value = 34417.5
$ws.Range("N88").Value = -35229.5
$ws.Range("H91").Value = 34417.5
$ws.Range("J91").Value = 34417.5
$ws.Range("L91").Value = 34417.5
$ws.Range("N91").Value = -37225.5
$ws.Range("H132").Value = 753054.2
$ws.Range("I132").Value = 1127745.8
$ws.Range("J132").Value = 3671
$ws.Range("K132").Value = 3383237.4
$ws.Range("L132").Value = 11013
$ws.Range("M132").Value = -3380707.4
$ws.Range("N132").Value = -16073
$ws.Range("H134").Value = 1204.1694
$ws.Range("I134").Value = 988.4666999999999
$ws.Range("J134").Value = 1897.5
$ws.Range("K134").Value = 2965.4001
$ws.Range("L134").Value = 5692.5
$ws.Range("M134").Value = -430.4000999999998
$ws.Range("N134").Value = -10762.5
$ws.Range("H136").Value = 757208.0600000001
$ws.Range("I136").Value = 950749
$ws.Range("J136").Value = 2398.3
$ws.Range("K136").Value = 2852247
$ws.Range("L136").Value = 7194.900000000001
$ws.Range("M136").Value = -2849697
$ws.Range("N136").Value = -12294.9

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 603901
$ws.Range("J12").Value = 805182.4399999999
$ws.Range("L12").Value = 2415547.32
$ws.Range("N12").Value = -2415893.32
$ws.Range("H75").Value = 10982.6
$ws.Range("I75").Value = 4956.5
$ws.Range("J75").Value = 15000
$ws.Range("K75").Value = 14869.5
$ws.Range("L75").Value = 45000
$ws.Range("M75").Value = -13871.5
$ws.Range("N75").Value = -46996
$ws.Range("H78").Value = 10982.6
$ws.Range("I78").Value = 4956.5
$ws.Range("J78").Value = 15000
$ws.Range("K78").Value = 44608.5
$ws.Range("L78").Value = 135000
$ws.Range("M78").Value = -39616.5
$ws.Range("N78").Value = -144984
$ws.Range("H92").Value = 514.8889
$ws.Range("I92").Value = 565.6
$ws.Range("K92").Value = 1696.8
$ws.Range("M92").Value = -448.8000000000002
$ws.Range("H130").Value = 2740
$ws.Range("I130").Value = 920
$ws.Range("K130").Value = 2760
$ws.Range("M130").Value = 2260

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 6671333.5
$ws.Range("I44").Value = 20000000
$ws.Range("J44").Value = 7000
$ws.Range("K44").Value = 20000000
$ws.Range("L44").Value = 7000
$ws.Range("M44").Value = -19999404
$ws.Range("N44").Value = -8192
$ws.Range("H70").Value = 6206.8965
$ws.Range("I70").Value = 5761.5
$ws.Range("J70").Value = 6935.727
$ws.Range("K70").Value = 5761.5
$ws.Range("L70").Value = 6935.727
$ws.Range("M70").Value = -5491.5
$ws.Range("N70").Value = -7475.727
$ws.Range("H73").Value = 6206.8965
$ws.Range("I73").Value = 5761.5
$ws.Range("J73").Value = 6935.727
$ws.Range("K73").Value = 5761.5
$ws.Range("L73").Value = 6935.727
$ws.Range("M73").Value = -4825.5
$ws.Range("N73").Value = -8807.726999999999
$ws.Range("H132").Value = 2008.3939
$ws.Range("I132").Value = 1399.3684
$ws.Range("K132").Value = 4198.1052
$ws.Range("M132").Value = -1668.1052

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 24147.428
$ws.Range("I30").Value = 24147.428
$ws.Range("K30").Value = 24147.428
$ws.Range("M30").Value = -24039.428
$ws.Range("H41").Value = 28333.334
$ws.Range("I41").Value = 28333.334
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 28333.334
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -27895.334
$ws.Range("N41").ClearContents()
$ws.Range("H132").Value = 4919.5
$ws.Range("I132").Value = 4667.3335
$ws.Range("J132").Value = 5676
$ws.Range("K132").Value = 14002.0005
$ws.Range("L132").Value = 17028
$ws.Range("M132").Value = -11472.0005
$ws.Range("N132").Value = -22088
$ws.Range("H136").Value = 2552.4856
$ws.Range("I136").Value = 1977.7333
$ws.Range("K136").Value = 5933.199900000001
$ws.Range("M136").Value = -3383.199900000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H42").Value = 39014.668
$ws.Range("J42").Value = 17000
$ws.Range("L42").Value = 17000
$ws.Range("N42").Value = -17756
$ws.Range("H43").Value = 19950
$ws.Range("I43").Value = 27900
$ws.Range("J43").Value = 12000
$ws.Range("K43").Value = 27900
$ws.Range("L43").Value = 12000
$ws.Range("M43").Value = -27751
$ws.Range("N43").Value = -12298
$ws.Range("H132").Value = 2238.95
$ws.Range("I132").Value = 1626.8572
$ws.Range("J132").Value = 3667.1667
$ws.Range("K132").Value = 4880.571599999999
$ws.Range("L132").Value = 11001.5001
$ws.Range("M132").Value = -2350.571599999999
$ws.Range("N132").Value = -16061.5001
$ws.Range("H136").Value = 1453.8
$ws.Range("I136").Value = 1453.8
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4361.4
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1811.4
$ws.Range("N136").ClearContents()
